$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row at row 38 (pushes the former rows 38-51 down
# to become rows 39-52, which is exactly what the target workbook shows).
$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with its values.
$ws.Cells.Item(38, 1).Value = 5
$ws.Cells.Item(38, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(38, 3).Value = "Maule"
$ws.Cells.Item(38, 4).Value = 44917
$ws.Cells.Item(38, 5).Value = 7
$ws.Cells.Item(38, 6).Value = "Fruta"
$ws.Cells.Item(38, 7).Value = 100103
$ws.Cells.Item(38, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(38, 9).Value = 100103003
$ws.Cells.Item(38, 10).Value = "Damasco"
$ws.Cells.Item(38, 11).Value = "Dina"
$ws.Cells.Item(38, 12).Value = "Primera"
$ws.Cells.Item(38, 13).Value = 200
$ws.Cells.Item(38, 14).Value = 17000
$ws.Cells.Item(38, 15).Value = 17000
$ws.Cells.Item(38, 16).Value = 17000
$ws.Cells.Item(38, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(38, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(38, 19).Value = 944
$ws.Cells.Item(38, 20).Value = 18
